$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.954.03"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.779.75"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.63"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.84"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.54%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.19%  "
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.94"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E11").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.35"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.60"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.224.17"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.784.61"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.922"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.865.68"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.14"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.27%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.64"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.58"
$ws.Range("D24").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.53"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E26").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E28").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E29").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.29"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.86"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.70"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.40"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.75%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0836"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E36").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.14"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.00"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.55"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.20"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.75"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.30%  "
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.064.68"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.70"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.924"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.93"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.51%  "
$ws.Range("E51").ClearFormats()
